$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Update raw Speed/Acceleration/Turning scores for several drivers (more realistic stats)
$ws.Range("D3").Value = 94
$ws.Range("E3").Value = 95
$ws.Range("F3").Value = 95

$ws.Range("D4").Value = 97
$ws.Range("E4").Value = 95
$ws.Range("F4").Value = 93

$ws.Range("D5").Value = 94
$ws.Range("E5").Value = 94
$ws.Range("F5").Value = 94

$ws.Range("D6").Value = 92
$ws.Range("E6").Value = 90
$ws.Range("F6").Value = 90

$ws.Range("D8").Value = 80
$ws.Range("E8").Value = 84
$ws.Range("F8").Value = 83

$ws.Range("D9").Value = 93
$ws.Range("E9").Value = 94
$ws.Range("F9").Value = 90

$ws.Range("D10").Value = 96
$ws.Range("E10").Value = 96
$ws.Range("F10").Value = 96

$ws.Range("D11").Value = 82
$ws.Range("E11").Value = 75
$ws.Range("F11").Value = 82

$ws.Range("D12").Value = 86
$ws.Range("E12").Value = 86
$ws.Range("F12").Value = 85

$ws.Range("D13").Value = 76
$ws.Range("E13").Value = 82
$ws.Range("F13").Value = 78

$ws.Range("D16").Value = 80
$ws.Range("E16").Value = 74
$ws.Range("F16").Value = 82

$ws.Range("D18").Value = 76
$ws.Range("E18").Value = 76
$ws.Range("F18").Value = 79

$ws.Range("D19").Value = 85
$ws.Range("E19").Value = 85
$ws.Range("F19").Value = 82

$ws.Range("D20").Value = 78
$ws.Range("E20").Value = 73
$ws.Range("F20").Value = 80

# Re-sort the table by Average Performance (descending) and turn on AutoFilter
$tableRange = $ws.Range("B2:H32")
$sortKey = $ws.Range("G2:G32")
$tableRange.Sort($sortKey, 2, $null, $null, 1, $null, $null, 1)
